$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PWM")

# --- Colors header -------------------------------------------------
$ws.Range("A12").Value = "Colors"
$ws.Range("A12").Style = "Accent1"

# --- Color maxima (Good style) --------------------------------------
$ws.Range("A13").Value = "RED_PWM_MAX"
$ws.Range("B13").Value = 160
$ws.Range("A13:B13").Style = "Good"

$ws.Range("A14").Value = "GREEN_PWM_MAX"
$ws.Range("B14").Value = 160
$ws.Range("A14:B14").Style = "Good"

$ws.Range("A15").Value = "BLUE_PWM_MAX"
$ws.Range("B15").Value = 160
$ws.Range("A15:B15").Style = "Good"

$ws.Range("A16").Value = "COLOR_STEP"
$ws.Range("B16").Value = 8
$ws.Range("A16:B16").Style = "Good"

# --- Variations (Neutral style) -------------------------------------
$ws.Range("A17").Value = "Red variations"
$ws.Range("B17").Formula = "=B13/`$B`$16"
$ws.Range("A17:B17").Style = "Neutral"

$ws.Range("A18").Value = "Green variations"
$ws.Range("A19").Value = "Blue variations"
$ws.Range("B18:B19").Formula = "=B14/`$B`$16"
$ws.Range("A18:B19").Style = "Neutral"

# --- Colors count (Calculation style) --------------------------------
$ws.Range("A20").Value = "Colors count"
$ws.Range("B20").Formula = "=B19*B18*B17"
$ws.Range("A20:B20").Style = "Calculation"

# --- Make PWM the active tab/sheet with B21 selected ------------------
$ws.Activate()
$ws.Range("B21").Select() | Out-Null
